# Updates cryptos list (coin name/link/price/volume) to match the
# latest scrape, per the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('E2').NumberFormat = "@"
$ws.Range('D2').Value = '22.021.31'
$ws.Range('E2').Value = '  -1.92%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.554.13'

# Row 4
$ws.Range('D4').NumberFormat = "@"
$ws.Range('E4').NumberFormat = "@"
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.00%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('E5').NumberFormat = "@"
$ws.Range('D5').Value = '1.000'
$ws.Range('E5').Value = '  -0.03%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('E6').NumberFormat = "@"
$ws.Range('D6').Value = '286.57'
$ws.Range('E6').Value = '  -0.45%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('E7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3758'
$ws.Range('E7').Value = '  +1.33%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('E8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3237'
$ws.Range('E8').Value = '  -2.42%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('E9').NumberFormat = "@"
$ws.Range('B9').Value = 'Polygon'
$ws.Range('C9').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D9').Value = '1.125'
$ws.Range('E9').Value = '  -2.38%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('E10').NumberFormat = "@"
$ws.Range('B10').Value = 'OKB'
$ws.Range('C10').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D10').Value = '41.09'
$ws.Range('E10').Value = '  -13.31%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('E11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07299'
$ws.Range('E11').Value = '  -2.70%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('E12').NumberFormat = "@"
$ws.Range('D12').Value = '1.001'
$ws.Range('E12').Value = '  +0.01%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('E13').NumberFormat = "@"
$ws.Range('D13').Value = '19.53'
$ws.Range('E13').Value = '  -5.94%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('E14').NumberFormat = "@"
$ws.Range('D14').Value = '5.702'
$ws.Range('E14').Value = '  -3.84%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('E15').NumberFormat = "@"
$ws.Range('D15').Value = '6.834'
$ws.Range('E15').Value = '  -1.25%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('E16').NumberFormat = "@"
$ws.Range('D16').Value = '1.554.29'
$ws.Range('E16').Value = '  -0.60%  '

# Row 17
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -3.08%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('E18').NumberFormat = "@"
$ws.Range('D18').Value = '0.06648'
$ws.Range('E18').Value = '  -1.19%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('E19').NumberFormat = "@"
$ws.Range('D19').Value = '85.11'
$ws.Range('E19').Value = '  -3.65%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('E20').NumberFormat = "@"
$ws.Range('D20').Value = '6.438'
$ws.Range('E20').Value = '  +0.88%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('E21').NumberFormat = "@"
$ws.Range('D21').Value = '1.0000'
$ws.Range('E21').Value = '  +0.01%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('E22').NumberFormat = "@"
$ws.Range('D22').Value = '15.97'
$ws.Range('E22').Value = '  -3.04%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('E23').NumberFormat = "@"
$ws.Range('D23').Value = '11.55'
$ws.Range('E23').Value = '  -3.58%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('E24').NumberFormat = "@"
$ws.Range('D24').Value = '22.046.47'
$ws.Range('E24').Value = '  -1.72%  '

# Row 25
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -6.07%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('E26').NumberFormat = "@"
$ws.Range('D26').Value = '2.524'
$ws.Range('E26').Value = '  -3.91%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('E27').NumberFormat = "@"
$ws.Range('D27').Value = '150.02'
$ws.Range('E27').Value = '  -0.24%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('E28').NumberFormat = "@"
$ws.Range('D28').Value = '18.85'
$ws.Range('E28').Value = '  -3.75%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('E29').NumberFormat = "@"
$ws.Range('D29').Value = '4.836'
$ws.Range('E29').Value = '  -2.38%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('E30').NumberFormat = "@"
$ws.Range('D30').Value = '1.729.98'
$ws.Range('E30').Value = '  -0.72%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('E31').NumberFormat = "@"
$ws.Range('D31').Value = '120.05'
$ws.Range('E31').Value = '  -4.03%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('E32').NumberFormat = "@"
$ws.Range('D32').Value = '1.121'
$ws.Range('E32').Value = '  +2.30%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('E33').NumberFormat = "@"
$ws.Range('D33').Value = '5.916'
$ws.Range('E33').Value = '  -2.73%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('E34').NumberFormat = "@"
$ws.Range('D34').Value = '1.708'
$ws.Range('E34').Value = '  -14.01%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('E35').NumberFormat = "@"
$ws.Range('D35').Value = '9.283'
$ws.Range('E35').Value = '  -5.94%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('E36').NumberFormat = "@"
$ws.Range('D36').Value = '0.08137'
$ws.Range('E36').Value = '  -2.31%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('E37').NumberFormat = "@"
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '0.02287'
$ws.Range('E37').Value = '  -6.43%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('E38').NumberFormat = "@"
$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').Value = '5.223'
$ws.Range('E38').Value = '  -1.83%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('E39').NumberFormat = "@"
$ws.Range('D39').Value = '0.06143'
$ws.Range('E39').Value = '  -3.72%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('E40').NumberFormat = "@"
$ws.Range('D40').Value = '0.2114'
$ws.Range('E40').Value = '  -4.93%  '

# Row 41
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -6.50%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('E42').NumberFormat = "@"
$ws.Range('D42').Value = '10.90'
$ws.Range('E42').Value = '  -4.16%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('E43').NumberFormat = "@"
$ws.Range('D43').Value = '1.0000'
$ws.Range('E43').Value = '  -0.03%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('E44').NumberFormat = "@"
$ws.Range('D44').Value = '0.5941'
$ws.Range('E44').Value = '  -4.74%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('E45').NumberFormat = "@"
$ws.Range('D45').Value = '13.49'
$ws.Range('E45').Value = '  -3.87%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('E46').NumberFormat = "@"
$ws.Range('D46').Value = '3.723'
$ws.Range('E46').Value = '  -1.29%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('E47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5739'
$ws.Range('E47').Value = '  -5.34%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('E48').NumberFormat = "@"
$ws.Range('D48').Value = '1.945'
$ws.Range('E48').Value = '  -4.83%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('E49').NumberFormat = "@"
$ws.Range('D49').Value = '119.91'
$ws.Range('E49').Value = '  -3.92%  '

# Row 50
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -4.25%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('E51').NumberFormat = "@"
$ws.Range('D51').Value = '0.06930'
$ws.Range('E51').Value = '  -3.72%  '
